$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Delete the "EUchild" sheet
# ---------------------------------------------------------------------
$euchild = $wb.Worksheets.Item("EUchild")
$euchild.Delete()

# ---------------------------------------------------------------------
# 2. Summary sheet - append rows 8:10
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A8").Value = 4
$summary.Range("B8").Value = "EU"
$summary.Range("A9").Value = 5
$summary.Range("B9").Value = "panel B for BE"
$summary.Range("A10").Value = 6
$summary.Range("B10").Value = "panel B for NL"

# ---------------------------------------------------------------------
# 3. BE sheet - update waves for rows 13:15, append row 16
# ---------------------------------------------------------------------
$be = $wb.Worksheets.Item("BE")
$be.Range("B13").Value = 5
$be.Range("B14").Value = 5
$be.Range("B15").Value = 5

$be.Range("A16").Value = "be"
$be.Range("B16").Value = 5
$be.Range("C16").Value = 0
$be.Range("D16").Value = 15
$be.Range("E16").Value = "B"
$be.Range("F16").Value = 7
$be.Range("G16").Value = 44249
$be.Range("G16").NumberFormat = $be.Range("G15").NumberFormat
$be.Range("H16").Value = "20_060765_BE2_Wave7_Final_v1_220221_IntClientUse"
$be.Range("I16").Formula = '=A16&"_"&"wk"&TEXT(D16,"00")&"_"&YEAR(G16)&TEXT(G16,"MM")&TEXT(G16,"DD")&"_p"&E16&"_wv"&TEXT(F16,"00")&""'

# ---------------------------------------------------------------------
# 4. NL sheet - update wave numbers for rows 10:11, append rows 12:13
# ---------------------------------------------------------------------
$nl = $wb.Worksheets.Item("NL")
$nl.Range("A10").Value = 6
$nl.Range("A11").Value = 6

$nl.Range("A12").Value = 6
$nl.Range("B12").Value = 0
$nl.Range("C12").Value = "nl"
$nl.Range("D12").Value = 11
$nl.Range("E12").Value = "B"
$nl.Range("F12").Value = 3
$nl.Range("G12").Value = 44229
$nl.Range("G12").NumberFormat = $nl.Range("G11").NumberFormat
$nl.Range("H12").Value = "20-090916_NL_Wave3_Final_v1_020221_IntClientUse"
$nl.Range("I12").Formula = '=C12&"_"&"wk"&TEXT(D12,"00")&"_"&YEAR(G12)&TEXT(G12,"MM")&TEXT(G12,"DD")&"_p"&E12&"_wv"&TEXT(F12,"00")&""'

$nl.Range("A13").Value = 6
$nl.Range("B13").Value = 0
$nl.Range("C13").Value = "nl"
$nl.Range("D13").Value = 12
$nl.Range("E13").Value = "B"
$nl.Range("F13").Value = 4
$nl.Range("G13").Value = 44239
$nl.Range("G13").NumberFormat = $nl.Range("G11").NumberFormat
$nl.Range("H13").Value = "20-090916_NL_Wave4_Final_v1_12022021_IntClientUse"
$nl.Range("I13").Formula = '=C13&"_"&"wk"&TEXT(D13,"00")&"_"&YEAR(G13)&TEXT(G13,"MM")&TEXT(G13,"DD")&"_p"&E13&"_wv"&TEXT(F13,"00")&""'

# ---------------------------------------------------------------------
# 5. G1 sheet - bump wave counts from 5 to 4 (rows 4:27, column B)
# ---------------------------------------------------------------------
$g1 = $wb.Worksheets.Item("G1")
for ($r = 4; $r -le 27; $r++) {
    $g1.Cells.Item($r, 2).Value = 4
}

# ---------------------------------------------------------------------
# 6. G2 sheet - bump wave counts from 5 to 4 (rows 2:8, column B)
# ---------------------------------------------------------------------
$g2 = $wb.Worksheets.Item("G2")
for ($r = 2; $r -le 8; $r++) {
    $g2.Cells.Item($r, 2).Value = 4
}
